$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.527.08"
$ws.Range("E2").Value = "'  -0.88%  "

$ws.Range("D3").Value = "'2.546.05"
$ws.Range("E3").Value = "'  +0.47%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "'  -0.03%  "

$ws.Range("D5").Value = "'308.87"
$ws.Range("E5").Value = "'  -2.59%  "

$ws.Range("D6").Value = "'97.42"
$ws.Range("E6").Value = "'  -0.47%  "

$ws.Range("D7").Value = "'0.572"
$ws.Range("E7").Value = "'  -0.49%  "

$ws.Range("E8").Value = "'  +0.00%  "

$ws.Range("E9").Value = "'  -1.04%  "

$ws.Range("D10").Value = "'35.45"
$ws.Range("E10").Value = "'  -1.53%  "

$ws.Range("E11").Value = "'  -0.90%  "

$ws.Range("E12").Value = "'  -2.97%  "

$ws.Range("E13").Value = "'  -2.11%  "

$ws.Range("D14").Value = "'2.935.27"
$ws.Range("E14").Value = "'  +0.59%  "

$ws.Range("D15").Value = "'15.80"
$ws.Range("E15").Value = "'  +4.31%  "

$ws.Range("D16").Value = "'2.577.69"
$ws.Range("E16").Value = "'  +1.95%  "

$ws.Range("D17").Value = "'0.835"
$ws.Range("E17").Value = "'  -1.82%  "

$ws.Range("D18").Value = "'42.622.72"
$ws.Range("E18").Value = "'  -0.70%  "

$ws.Range("D19").Value = "'6.74"
$ws.Range("E19").Value = "'  -2.09%  "

$ws.Range("B20").Value = "'InternetComputer(DFINITY)"
$ws.Range("C20").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.37"
$ws.Range("E20").Value = "'  -2.90%  "

$ws.Range("B21").Value = "'ShibaInu"
$ws.Range("C21").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0₃0955"
$ws.Range("E21").Value = "'  -0.91%  "

$ws.Range("D22").Value = "'69.27"

$ws.Range("D23").Value = "'247.29"
$ws.Range("E23").Value = "'  -2.19%  "

$ws.Range("D24").Value = "'2.91"
$ws.Range("E24").Value = "'  -1.65%  "

$ws.Range("E25").Value = "'  +0.05%  "

$ws.Range("D26").Value = "'26.62"
$ws.Range("E26").Value = "'  +0.57%  "

$ws.Range("E27").Value = "'  +0.03%  "

$ws.Range("D28").Value = "'2.36"
$ws.Range("E28").Value = "'  -1.69%  "

$ws.Range("D29").Value = "'40.39"
$ws.Range("E29").Value = "'  -1.80%  "

$ws.Range("D30").Value = "'10.11"
$ws.Range("E30").Value = "'  -3.03%  "

$ws.Range("D31").Value = "'157.99"
$ws.Range("E31").Value = "'  +0.20%  "

$ws.Range("D32").Value = "'5.72"
$ws.Range("E32").Value = "'  -3.49%  "

$ws.Range("D33").Value = "'0.0794"
$ws.Range("E33").Value = "'  +0.66%  "

$ws.Range("E34").Value = "'  -1.84%  "

$ws.Range("D35").Value = "'2.08"
$ws.Range("E35").Value = "'  -3.87%  "

$ws.Range("D36").Value = "'2.62"
$ws.Range("E36").Value = "'  -3.32%  "

$ws.Range("B37").Value = "'ApeXProtocol"
$ws.Range("C37").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D37").Value = "'2.58"
$ws.Range("E37").Value = "'  +4.39%  "

$ws.Range("B38").Value = "'Celestia"
$ws.Range("C38").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'18.30"
$ws.Range("E38").Value = "'  -3.78%  "

$ws.Range("E39").Value = "'  -1.80%  "

$ws.Range("D40").Value = "'0.118"
$ws.Range("E40").Value = "'  -0.63%  "

$ws.Range("D41").Value = "'22.39"
$ws.Range("E41").Value = "'  +2.68%  "

$ws.Range("D42").Value = "'4.05"
$ws.Range("E42").Value = "'  +5.34%  "

$ws.Range("E43").Value = "'  -0.17%  "

$ws.Range("D44").Value = "'0.0299"
$ws.Range("E44").Value = "'  -1.70%  "

$ws.Range("D45").Value = "'1.991.23"
$ws.Range("E45").Value = "'  -1.34%  "

$ws.Range("E46").Value = "'  -3.21%  "

$ws.Range("D47").Value = "'9.02"
$ws.Range("E47").Value = "'  -0.76%  "

$ws.Range("D48").Value = "'2.789.26"
$ws.Range("E48").Value = "'  +0.52%  "

$ws.Range("D49").Value = "'80.90"
$ws.Range("E49").Value = "'  -4.13%  "

$ws.Range("E50").Value = "'  -0.49%  "

$ws.Range("D51").Value = "'73.32"
$ws.Range("E51").Value = "'  -4.34%  "
